# Apply the "files updated and bug fixed" change to the Gori sheet:
# update the row-4 figures for years 2015-2021 (columns E:K) and
# leave the selection on the newly edited range (E4:K4, active cell E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gori")

$newValues = @{
    "E4" = 21340
    "F4" = 21674
    "G4" = 22026
    "H4" = 22303
    "I4" = 22793
    "J4" = 23545
    "K4" = 24002
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

$ws.Activate()
$ws.Range("E4:K4").Select()
